# Generate Report for Handoff
#
# For the file "a8467e36-9db8-48dc-9000-7fa5ce1d8d85" the handoff transform
# was re-run, producing a fresh "Latest Handoff Datetime" timestamp. That
# file is a dependency of several other rows, so its newly generated handoff
# timestamp is reflected in the "Latest Handoff Datetime" column (column D)
# for all of the rows that depend on it, on both the "zh-cn" and "de-de"
# report sheets.

$wb = $excel.ActiveWorkbook

$rows = @(4, 6, 7, 8, 9, 10)

$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Cells.Item($r, 4).Value = "2016-02-25 06:33:41"
}

$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Cells.Item($r, 4).Value = "2016-02-25 06:33:54"
}
